# Add a new hero-skill row ("skill9" / "亡灵" undead summon) to the
# HeroSkill table, growing the table from A1:F11 to A1:F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing table by one row - this expands the table range,
# the autofilter range and the sheet dimension automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Fill the new row (row 12). Write column F (Icon) before the Chinese
# text columns so the new shared strings land in the order:
# skill9, 亡灵, 在一个坟墓上召唤一个骷髅.
$ws.Cells.Item(12, 1).Value = 31000009
$ws.Cells.Item(12, 6).Value = "skill9"
$ws.Cells.Item(12, 2).Value = "亡灵"
$ws.Cells.Item(12, 3).Value = "在一个坟墓上召唤一个骷髅"
$ws.Cells.Item(12, 4).Value = 3
$ws.Cells.Item(12, 5).Value = 53100006

# Match the saved selection state (active cell E12).
$ws.Range("E12").Select()
